$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 87, shifting existing rows 87-106 down to 88-107
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new weekly record
$ws.Cells.Item(87, 1).Value2 = 4
$ws.Cells.Item(87, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(87, 3).Value2 = "Los Lagos"
$ws.Cells.Item(87, 4).Value2 = 44637
$ws.Cells.Item(87, 5).Value2 = 10
$ws.Cells.Item(87, 6).Value2 = 100112052
$ws.Cells.Item(87, 7).Value2 = "Albahaca"
$ws.Cells.Item(87, 8).Value2 = "Sin especificar"
$ws.Cells.Item(87, 9).Value2 = "Primera"
$ws.Cells.Item(87, 10).Value2 = 90
$ws.Cells.Item(87, 11).Value2 = 7000
$ws.Cells.Item(87, 12).Value2 = 7000
$ws.Cells.Item(87, 13).Value2 = 7000
$ws.Cells.Item(87, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(87, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(87, 16).Value2 = 1167
$ws.Cells.Item(87, 17).Value2 = 6
$ws.Cells.Item(87, 18).Value2 = "Hortaliza"
